# Insert 2 new rows at position 1145, shifting existing rows 1145-1225 down to 1147-1227.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1145:1146").Insert()

# Row 1145 - new record (same constant columns as the rest of this dataset)
$ws.Cells.Item(1145, 1).Value = 9
$ws.Cells.Item(1145, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(1145, 3).Value = "Metropolitana"
$ws.Cells.Item(1145, 4).Value = 44746
$ws.Cells.Item(1145, 5).Value = 13
$ws.Cells.Item(1145, 6).Value = 100114001
$ws.Cells.Item(1145, 7).Value = "Papa"
$ws.Cells.Item(1145, 8).Value = "Asterix"
$ws.Cells.Item(1145, 9).Value = "1a (guarda lavada)"
$ws.Cells.Item(1145, 10).Value = 250
$ws.Cells.Item(1145, 11).Value = 9000
$ws.Cells.Item(1145, 12).Value = 10000
$ws.Cells.Item(1145, 13).Value = 9500
$ws.Cells.Item(1145, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(1145, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(1145, 16).Value = 380
$ws.Cells.Item(1145, 17).Value = 25
$ws.Cells.Item(1145, 18).Value = "Hortaliza"

# Row 1146 - new record
$ws.Cells.Item(1146, 1).Value = 9
$ws.Cells.Item(1146, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(1146, 3).Value = "Metropolitana"
$ws.Cells.Item(1146, 4).Value = 44746
$ws.Cells.Item(1146, 5).Value = 13
$ws.Cells.Item(1146, 6).Value = 100114001
$ws.Cells.Item(1146, 7).Value = "Papa"
$ws.Cells.Item(1146, 8).Value = "Asterix"
$ws.Cells.Item(1146, 9).Value = "1a (guarda)"
$ws.Cells.Item(1146, 10).Value = 160
$ws.Cells.Item(1146, 11).Value = 8000
$ws.Cells.Item(1146, 12).Value = 8000
$ws.Cells.Item(1146, 13).Value = 8000
$ws.Cells.Item(1146, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(1146, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(1146, 16).Value = 320
$ws.Cells.Item(1146, 17).Value = 25
$ws.Cells.Item(1146, 18).Value = "Hortaliza"
